# "clean up of atf code"
# Updates the protocoltestcasedetails sheet:
#  - testcase28 renamed from "csv_sample1_csv_sample2" to "manual_sql_etltesting"
#  - testcase29 / testcase30 rows removed (rows 30 & 31 cleared)
#  - two trailing blank rows (38 & 39) deleted, shrinking the used range
#  - the Y/N data-validation range shrunk to match the new last data row
#  - selection/scroll position updated

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protocoltestcasedetails")

# Rename testcase28's test case name; the dependent CONCAT formula in column C
# recalculates automatically.
$ws.Range("B29").Value = "testcase28_manual_sql_etltesting"

# testcase29 and testcase30 rows (30 & 31) are no longer used - clear them out
# but keep the formatting/style already applied to those cells.
$ws.Range("A30:D31").ClearContents()

# Remove the two now-superfluous trailing blank rows (38 & 39) entirely so the
# sheet's used range shrinks from D39 to D37.
$ws.Range("A38:D39").Delete()

# The list validation on column D previously covered D2:D31; shrink it down to
# the new last populated row D29.
$ws.Range("D2:D31").Validation.Delete()
$ws.Range("D2:D29").Validation.Add(3, 1, 1, """Y,N""")

# Update the view so the previously-used row 33 area is visible/selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B33").Select()
